$d = $word.ActiveDocument

# The document has two "Pearson Edexcel" logo pictures (in the footers,
# currently saved internally as image2.png) and two "BTec_Logo-Orange"
# pictures (in the headers, currently saved internally as image1.jpg).
# This rename swaps each pair's internal picture name:
#   Pearson logo (footers): image2.png -> image1.png
#   BTec logo   (headers) : image1.jpg -> image2.jpg

function Rename-LogoShape($shape) {
    $desc = $shape.AlternativeText
    if ($desc -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shape.Name = "image1.png"
    } elseif ($desc -eq "BTec_Logo-Orange") {
        $shape.Name = "image2.jpg"
    }
}

foreach ($sec in $d.Sections) {
    for ($hidx = 1; $hidx -le 3; $hidx++) {
        $hf = $sec.Headers.Item($hidx)
        if ($hf.Exists) {
            $rng = $hf.Range
            for ($i = 1; $i -le $rng.InlineShapes.Count; $i++) {
                Rename-LogoShape $rng.InlineShapes.Item($i)
            }
        }
    }
    for ($fidx = 1; $fidx -le 3; $fidx++) {
        $ft = $sec.Footers.Item($fidx)
        if ($ft.Exists) {
            $rng = $ft.Range
            for ($i = 1; $i -le $rng.InlineShapes.Count; $i++) {
                Rename-LogoShape $rng.InlineShapes.Item($i)
            }
        }
    }
}
